# Auto update: 2025-12-05 03:03:39
#
# Daily refresh of the US-market quantum-computing decision table on
# Sheet1. The scoring run re-ranked the watch list (Rigetti now scores
# above D-Wave and IBM moves to the bottom) and refreshed the recomputed
# metric columns (종가/RSI/5일수익률/최종점수 etc.) for every ticker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: IonQ, Inc. (IONQ) stays on top; metrics refreshed ---
$ws.Range("D2").Value = 54.79
$ws.Range("E2").Value = 65.3
$ws.Range("F2").Value = 16.82
$ws.Range("N2").Value = 53.71147335634279

# --- Row 3: now Rigetti Computing, Inc. (RGTI) ---
$ws.Range("B3").Value = "Rigetti Computing, Inc."
$ws.Range("C3").Value = "RGTI"
$ws.Range("D3").Value = 29.18
$ws.Range("E3").Value = 61.3
$ws.Range("F3").Value = 14.14
$ws.Range("I3").Value = 63
$ws.Range("J3").Value = 83
$ws.Range("K3").Value = 56.3
$ws.Range("N3").Value = 53.71147335634279

# --- Row 4: now D-Wave Quantum Inc. (QBTS) ---
$ws.Range("B4").Value = "D-Wave Quantum Inc."
$ws.Range("C4").Value = "QBTS"
$ws.Range("D4").Value = 28.33
$ws.Range("E4").Value = 65
$ws.Range("F4").Value = 26.44
$ws.Range("H4").Value = 70
$ws.Range("J4").Value = 76
$ws.Range("N4").Value = 53.71147335634279

# --- Row 5: now International Business Machines (IBM) ---
$ws.Range("B5").Value = "International Business Machines"
$ws.Range("C5").Value = "IBM"
$ws.Range("D5").Value = 306.73
$ws.Range("E5").Value = 51.8
$ws.Range("F5").Value = 1.16
$ws.Range("H5").Value = 66
$ws.Range("I5").Value = 63
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 56.3
$ws.Range("N5").Value = 53.71147335634279
